# Apply cryptos list price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Price" column (D) ---
# Column D values look numeric (contain dots used as thousands separators
# or plain decimals) but must stay stored as literal text, exactly as in the
# source data. Temporarily mark the cells as Text ("@") before assigning the
# value so Excel does not reinterpret/round them as numbers, then restore the
# original (default) cell style so formatting stays untouched.
$priceCells = @(
    @{ Ref = "D2"; Value = '27.461.57' }
    @{ Ref = "D3"; Value = '1.837.88' }
    @{ Ref = "D5"; Value = '314.28' }
    @{ Ref = "D6"; Value = '1.010' }
    @{ Ref = "D7"; Value = '0.4733' }
    @{ Ref = "D8"; Value = '0.3694' }
    @{ Ref = "D9"; Value = '0.07465' }
    @{ Ref = "D10"; Value = '0.8856' }
    @{ Ref = "D11"; Value = '20.45' }
    @{ Ref = "D12"; Value = '1.887.53' }
    @{ Ref = "D13"; Value = '0.07342' }
    @{ Ref = "D14"; Value = '5.451' }
    @{ Ref = "D15"; Value = '93.30' }
    @{ Ref = "D16"; Value = '6.583' }
    @{ Ref = "D17"; Value = '1.011' }
    @{ Ref = "D18"; Value = '0.000008825' }
    @{ Ref = "D19"; Value = '1.010' }
    @{ Ref = "D20"; Value = '14.80' }
    @{ Ref = "D21"; Value = '27.474.97' }
    @{ Ref = "D22"; Value = '5.325' }
    @{ Ref = "D24"; Value = '2.134.19' }
    @{ Ref = "D25"; Value = '1.907' }
    @{ Ref = "D26"; Value = '152.07' }
    @{ Ref = "D27"; Value = '18.63' }
    @{ Ref = "D28"; Value = '2.151' }
    @{ Ref = "D29"; Value = '5.251' }
    @{ Ref = "D30"; Value = '117.72' }
    @{ Ref = "D31"; Value = '0.08966' }
    @{ Ref = "D32"; Value = '0.7570' }
    @{ Ref = "D34"; Value = '4.556' }
    @{ Ref = "D35"; Value = '2.942' }
    @{ Ref = "D36"; Value = '1.012' }
    @{ Ref = "D37"; Value = '1.104' }
    @{ Ref = "D38"; Value = '0.05339' }
    @{ Ref = "D40"; Value = '2.992' }
    @{ Ref = "D41"; Value = '7.333' }
    @{ Ref = "D42"; Value = '2.408' }
    @{ Ref = "D44"; Value = '0.1662' }
    @{ Ref = "D45"; Value = '8.518' }
    @{ Ref = "D46"; Value = '0.4913' }
    @{ Ref = "D47"; Value = '10.53' }
    @{ Ref = "D49"; Value = '105.04' }
    @{ Ref = "D51"; Value = '0.06308' }
)
foreach ($item in $priceCells) {
    $cell = $ws.Range($item.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}

# --- Update "Volume(1h)" column (E) ---
$volumeCells = @(
    @{ Ref = "E2"; Value = '  +2.05%  ' }
    @{ Ref = "E4"; Value = '  +1.08%  ' }
    @{ Ref = "E5"; Value = '  +1.61%  ' }
    @{ Ref = "E6"; Value = '  +0.90%  ' }
    @{ Ref = "E7"; Value = '  +1.74%  ' }
    @{ Ref = "E8"; Value = '  +0.84%  ' }
    @{ Ref = "E9"; Value = '  +1.60%  ' }
    @{ Ref = "E10"; Value = '  +1.92%  ' }
    @{ Ref = "E11"; Value = '  +0.63%  ' }
    @{ Ref = "E12"; Value = '  +1.38%  ' }
    @{ Ref = "E13"; Value = '  +3.55%  ' }
    @{ Ref = "E14"; Value = '  +1.24%  ' }
    @{ Ref = "E15"; Value = '  +1.81%  ' }
    @{ Ref = "E16"; Value = '  +1.07%  ' }
    @{ Ref = "E17"; Value = '  +0.85%  ' }
    @{ Ref = "E18"; Value = '  +1.25%  ' }
    @{ Ref = "E19"; Value = '  +0.88%  ' }
    @{ Ref = "E20"; Value = '  +0.95%  ' }
    @{ Ref = "E21"; Value = '  +1.99%  ' }
    @{ Ref = "E22"; Value = '  +0.45%  ' }
    @{ Ref = "E23"; Value = '  +0.46%  ' }
    @{ Ref = "E24"; Value = '  +2.48%  ' }
    @{ Ref = "E25"; Value = '  +0.75%  ' }
    @{ Ref = "E26"; Value = '  +0.71%  ' }
    @{ Ref = "E27"; Value = '  +1.77%  ' }
    @{ Ref = "E28"; Value = '  +0.51%  ' }
    @{ Ref = "E29"; Value = '  -0.11%  ' }
    @{ Ref = "E30"; Value = '  +2.02%  ' }
    @{ Ref = "E31"; Value = '  +0.53%  ' }
    @{ Ref = "E32"; Value = '  -0.01%  ' }
    @{ Ref = "E33"; Value = '  +2.11%  ' }
    @{ Ref = "E34"; Value = '  +1.54%  ' }
    @{ Ref = "E35"; Value = '  +1.09%  ' }
    @{ Ref = "E36"; Value = '  +1.08%  ' }
    @{ Ref = "E37"; Value = '  +1.83%  ' }
    @{ Ref = "E38"; Value = '  +1.19%  ' }
    @{ Ref = "E39"; Value = '  +0.47%  ' }
    @{ Ref = "E40"; Value = '  +0.42%  ' }
    @{ Ref = "E41"; Value = '  +1.23%  ' }
    @{ Ref = "E42"; Value = '  +4.99%  ' }
    @{ Ref = "E43"; Value = '  +0.62%  ' }
    @{ Ref = "E44"; Value = '  +0.58%  ' }
    @{ Ref = "E45"; Value = '  +1.00%  ' }
    @{ Ref = "E46"; Value = '  +0.70%  ' }
    @{ Ref = "E47"; Value = '  +1.48%  ' }
    @{ Ref = "E48"; Value = '  +1.04%  ' }
    @{ Ref = "E49"; Value = '  +1.68%  ' }
    @{ Ref = "E50"; Value = '  +1.13%  ' }
    @{ Ref = "E51"; Value = '  +0.28%  ' }
)
foreach ($item in $volumeCells) {
    $ws.Range($item.Ref).Value = $item.Value
}
